# Weekly update: a new Haba price observation (date 2021-11-16) is inserted
# before the existing row 91, pushing the subsequent rows (old 91-94) down
# to 92-95 and extending the used range to A1:R95.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 91, shifting existing rows 91-94 down to 92-95
$ws.Rows("91").Insert()

# Populate the newly inserted row 91 with the new weekly price observation
$ws.Range("A91").Value = 3
$ws.Range("B91").Value = "Femacal de La Calera"
$ws.Range("C91").Value = "Coquimbo"
$ws.Range("D91").Value = 44516
$ws.Range("E91").Value = 5
$ws.Range("F91").Value = 100112026
$ws.Range("G91").Value = "Haba"
$ws.Range("H91").Value = "Sin especificar"
$ws.Range("I91").Value = "Primera"
$ws.Range("J91").Value = 90
$ws.Range("K91").Value = 7000
$ws.Range("L91").Value = 7500
$ws.Range("M91").Value = 7278
$ws.Range("N91").Value = "$/saco 25 kilos"
$ws.Range("O91").Value = "Provincia de Quillota"
$ws.Range("P91").Value = 291
$ws.Range("Q91").Value = 25
$ws.Range("R91").Value = "Hortaliza"
